$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Flattened row-major replacement values (20 rows x 5 cols = 100 arithmetic
# drill cells) matching the commit "Update master to output generated at
# 60844e3".
$values = @(
  "19-16=",
  "34+15=",
  "34+12=",
  "60+39=",
  "12+69=",
  "93-25=",
  "18+23=",
  "1+93=",
  "2+50=",
  "40-22=",
  "2+68=",
  "45-24=",
  "52+1=",
  "21+5=",
  "73-62=",
  "73-66=",
  "90-8=",
  "9+46=",
  "30-15=",
  "63+11=",
  "82-3=",
  "28+58=",
  "96-71=",
  "23+44=",
  "4+19=",
  "9+70=",
  "57-16=",
  "64-9=",
  "0+4=",
  "35+47=",
  "0+8=",
  "41+3=",
  "78-4=",
  "39-13=",
  "62+20=",
  "2+92=",
  "32+23=",
  "86+0=",
  "4+61=",
  "34+23=",
  "39-2=",
  "8+54=",
  "44+27=",
  "26-21=",
  "96-72=",
  "53+39=",
  "21+61=",
  "42-24=",
  "62-37=",
  "38-38=",
  "62+15=",
  "1+88=",
  "12-4=",
  "97-63=",
  "89-24=",
  "81-54=",
  "40-7=",
  "22+60=",
  "62-51=",
  "63-22=",
  "67-1=",
  "19+3=",
  "26+37=",
  "21-5=",
  "96-69=",
  "10+8=",
  "36+63=",
  "53-10=",
  "97-4=",
  "37-16=",
  "74-9=",
  "72+26=",
  "18+24=",
  "16+67=",
  "55+17=",
  "50-27=",
  "39-25=",
  "86-73=",
  "10+64=",
  "85-33=",
  "44-26=",
  "5+72=",
  "25-9=",
  "50-20=",
  "97-1=",
  "84-56=",
  "74-51=",
  "45-13=",
  "2+2=",
  "89-64=",
  "48+51=",
  "56-10=",
  "71+6=",
  "33+14=",
  "24+59=",
  "30+69=",
  "35+43=",
  "47-32=",
  "65-8=",
  "25+60="
)

$cols = $t.Columns.Count
$rows = $t.Rows.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $t.Cell($r, $c).Range.Text = $values[$idx]
    $idx = $idx + 1
  }
}
Write-Output ("Updated " + $idx + " cells across " + $rows + " rows x " + $cols + " cols")